$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current column C ("Edges" header),
# pushing "Edges" and the data after it two columns to the right.
# This also carries the quote-prefixed empty cell that used to be at F3
# over to H3 automatically.
$null = $ws.Range("C:D").Insert()

# The insert shifted the old G2 value into I2 - clear that leftover cell,
# new data below fills the H column instead.
$null = $ws.Range("I2").ClearContents()

# New column headers for the freshly inserted columns.
$ws.Range("C1").Value = "Yc"
$ws.Range("D1").Value = "Zc"

# Row 2 (existing section) - fill the new Yc/Zc values and the rest of the row.
$ws.Range("B2").Value = 0.25
$ws.Range("C2").Value = 0.125
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 8
$ws.Range("H2").Value = 7

# Row 3 - second closed section's data.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0.25
$ws.Range("C3").Value = -0.125
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 6

# H3 needs to keep its inherited quote-prefix style while holding the
# numeric value 8. Setting .Value directly clears that style, so stash the
# style on a scratch cell, assign the value, then paste the style back.
$null = $ws.Range("H3").Copy($ws.Range("Z1"))
$ws.Range("H3").Value = 8
$null = $ws.Range("Z1").Copy()
$null = $ws.Range("H3").PasteSpecial(-4122)
$null = $ws.Range("Z1").Clear()

# Resize the new "Yc" column to fit its (narrower) content, same as the
# "Edges" column already does for its own content.
$null = $ws.Columns("C:C").AutoFit()

$null = $ws.Range("M6").Select()
